$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -6.007399999999999
$ws.Range("D8").Value = -8.817099999999989
$ws.Range("B12").Value = 5.522899999999998
$ws.Range("D12").Value = -7.758199999999995
$ws.Range("D14").Value = -8.635799999999998
$ws.Range("D22").Value = -7.758799999999995
